# upd fx_predict.xlsx by request
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data updates (sheet "Предсказание") ---
# Row 2 - USD/RUB
$ws.Range("B2").Value = 87
$ws.Range("C2").Value = 88

# Row 3 - EUR/USD
$ws.Range("B3").Value = 1.09

# Row 4 - EUR/RUB
$ws.Range("B4").Value = 94.4
$ws.Range("C4").Value = 94.16

# Row 5 - USD/CNY
$ws.Range("B5").Value = 7.2
$ws.Range("C5").Value = 7.2

# Row 6 - CNY/RUB
$ws.Range("B6").Value = 12.8
$ws.Range("C6").Value = 12.22

# Row 7 - USD/INR
$ws.Range("B7").Value = 83.3
$ws.Range("F7").Value = 82

# Row 8 - INR/RUB
$ws.Range("B8").Value = 1.04
$ws.Range("C8").Value = 1.06
$ws.Range("F8").Value = 1.1

# Row 9 - USD/TRY
$ws.Range("B9").Value = 28.7
$ws.Range("C9").Value = 29
$ws.Range("D9").Value = 30
$ws.Range("E9").Value = 31
$ws.Range("F9").Value = 32

# Row 10 - TRY/RUB
$ws.Range("B10").Value = 3.03
$ws.Range("C10").Value = 3.03
$ws.Range("D10").Value = 3
$ws.Range("E10").Value = 2.9
$ws.Range("F10").Value = 2.81

# Row 11 - USD/KZT
$ws.Range("B11").Value = 465
$ws.Range("C11").Value = 465
$ws.Range("D11").Value = 470
$ws.Range("E11").Value = 470
$ws.Range("F11").Value = 470

# Row 12 - KZT/RUB*100
$ws.Range("B12").Value = 18.71
$ws.Range("C12").Value = 18.920000000000002
$ws.Range("D12").Value = 19.149999999999999
$ws.Range("E12").Value = 19.149999999999999
$ws.Range("F12").Value = 19.149999999999999

# --- Column widths (A widened for longer label, B:F narrowed to fit the new figures) ---
$ws.Range("A1").EntireColumn.ColumnWidth = 17.3
$ws.Range("B1").EntireColumn.ColumnWidth = 5.65
$ws.Range("C1").EntireColumn.ColumnWidth = 5.65
$ws.Range("D1").EntireColumn.ColumnWidth = 5.65
$ws.Range("E1").EntireColumn.ColumnWidth = 5.65
$ws.Range("F1").EntireColumn.ColumnWidth = 5.65

# --- Selection moved to I5 ---
$ws.Range("I5").Select()
